$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.364.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.775.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.94"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4234"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3599"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07181"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8366"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.43"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.756.71"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.443"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.239"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06898"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.03"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008652"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.96"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.368.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.083"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.986.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.56"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.789"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -8.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.082"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.34"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.834"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +11.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08842"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7266"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.122"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.315"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.0000"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.733"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05108"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1610"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4914"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.603"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.325"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.47%  "
$ws.Range("B44").Value = "PaxosStandard"
$ws.Range("C44").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -30.80%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.021"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.57"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.15"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.632"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06175"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4439"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.95%  "
